$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - DecisionTreeClassifier
$ws.Range("C2").Value = 0.9552
$ws.Range("E2").Value = 0.9488
$ws.Range("G2").Value = 0.9623
$ws.Range("I2").Value = 0.9555

# Row 3 - RandomForestClassifier
$ws.Range("C3").Value = 0.9498
$ws.Range("E3").Value = 0.9647
$ws.Range("G3").Value = 0.9338
$ws.Range("I3").Value = 0.949

# Row 4 - LogisticRegression
$ws.Range("B4").Value = 0.8494
$ws.Range("C4").Value = 0.837
$ws.Range("D4").Value = 0.8957000000000001
$ws.Range("E4").Value = 0.8799
$ws.Range("F4").Value = 0.7907999999999999
$ws.Range("G4").Value = 0.7806
$ws.Range("H4").Value = 0.84
$ws.Range("I4").Value = 0.8273

# Row 5 - KNeighborsClassifier
$ws.Range("B5").Value = 0.8831
$ws.Range("C5").Value = 0.8337
$ws.Range("D5").Value = 0.9256
$ws.Range("E5").Value = 0.8767
$ws.Range("F5").Value = 0.8333
$ws.Range("G5").Value = 0.7766
$ws.Range("H5").Value = 0.877
$ws.Range("I5").Value = 0.8236

# Row 6 - New: Sequential (TensorFlow dense neural network model)
$ws.Range("A6").Value = "Sequential"
$ws.Range("B6").Value = 0.822
$ws.Range("C6").Value = 0.8257
$ws.Range("D6").Value = 0.8653
$ws.Range("E6").Value = 0.8695000000000001
$ws.Range("F6").Value = 0.7628
$ws.Range("G6").Value = 0.7664
$ws.Range("H6").Value = 0.8108
$ws.Range("I6").Value = 0.8147
